$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "55-47="
$t.Cell(1, 2).Range.Text = "49+35="
$t.Cell(1, 3).Range.Text = "66-29="
$t.Cell(1, 4).Range.Text = "96-38="
$t.Cell(1, 5).Range.Text = "57+16="
$t.Cell(2, 1).Range.Text = "70-41="
$t.Cell(2, 2).Range.Text = "72-28="
$t.Cell(2, 3).Range.Text = "58+9="
$t.Cell(2, 4).Range.Text = "54-36="
$t.Cell(2, 5).Range.Text = "45+26="
$t.Cell(3, 1).Range.Text = "50-49="
$t.Cell(3, 2).Range.Text = "15+69="
$t.Cell(3, 3).Range.Text = "57+34="
$t.Cell(3, 4).Range.Text = "40-16="
$t.Cell(3, 5).Range.Text = "79+17="
$t.Cell(4, 1).Range.Text = "48+4="
$t.Cell(4, 2).Range.Text = "45-38="
$t.Cell(4, 3).Range.Text = "44-37="
$t.Cell(4, 4).Range.Text = "80-18="
$t.Cell(4, 5).Range.Text = "53-5="
$t.Cell(5, 1).Range.Text = "56-28="
$t.Cell(5, 2).Range.Text = "79+6="
$t.Cell(5, 3).Range.Text = "56+39="
$t.Cell(5, 4).Range.Text = "90-15="
$t.Cell(5, 5).Range.Text = "43+9="
$t.Cell(6, 1).Range.Text = "17+15="
$t.Cell(6, 2).Range.Text = "67+15="
$t.Cell(6, 3).Range.Text = "16+25="
$t.Cell(6, 4).Range.Text = "73-6="
$t.Cell(6, 5).Range.Text = "90-13="
$t.Cell(7, 1).Range.Text = "48+45="
$t.Cell(7, 2).Range.Text = "47-29="
$t.Cell(7, 3).Range.Text = "49+44="
$t.Cell(7, 4).Range.Text = "92-87="
$t.Cell(7, 5).Range.Text = "49+46="
$t.Cell(8, 1).Range.Text = "36-17="
$t.Cell(8, 2).Range.Text = "19+78="
$t.Cell(8, 3).Range.Text = "71-36="
$t.Cell(8, 4).Range.Text = "60-27="
$t.Cell(8, 5).Range.Text = "15+38="
$t.Cell(9, 1).Range.Text = "16+45="
$t.Cell(9, 2).Range.Text = "74-55="
$t.Cell(9, 3).Range.Text = "95-77="
$t.Cell(9, 4).Range.Text = "19+7="
$t.Cell(9, 5).Range.Text = "49+8="
$t.Cell(10, 1).Range.Text = "27+9="
$t.Cell(10, 2).Range.Text = "29+68="
$t.Cell(10, 3).Range.Text = "36+6="
$t.Cell(10, 4).Range.Text = "24+8="
$t.Cell(10, 5).Range.Text = "84-47="
$t.Cell(11, 1).Range.Text = "81-76="
$t.Cell(11, 2).Range.Text = "60-17="
$t.Cell(11, 3).Range.Text = "93-77="
$t.Cell(11, 4).Range.Text = "9+8="
$t.Cell(11, 5).Range.Text = "19+57="
$t.Cell(12, 1).Range.Text = "51-43="
$t.Cell(12, 2).Range.Text = "67-19="
$t.Cell(12, 3).Range.Text = "8+27="
$t.Cell(12, 4).Range.Text = "90-39="
$t.Cell(12, 5).Range.Text = "25-9="
$t.Cell(13, 1).Range.Text = "27+29="
$t.Cell(13, 2).Range.Text = "56-8="
$t.Cell(13, 3).Range.Text = "52+29="
$t.Cell(13, 4).Range.Text = "9+77="
$t.Cell(13, 5).Range.Text = "85-7="
$t.Cell(14, 1).Range.Text = "8+45="
$t.Cell(14, 2).Range.Text = "37+57="
$t.Cell(14, 3).Range.Text = "19+15="
$t.Cell(14, 4).Range.Text = "73-19="
$t.Cell(14, 5).Range.Text = "17+26="
$t.Cell(15, 1).Range.Text = "7+66="
$t.Cell(15, 2).Range.Text = "19+34="
$t.Cell(15, 3).Range.Text = "39+52="
$t.Cell(15, 4).Range.Text = "18+13="
$t.Cell(15, 5).Range.Text = "81-78="
$t.Cell(16, 1).Range.Text = "36-7="
$t.Cell(16, 2).Range.Text = "45+28="
$t.Cell(16, 3).Range.Text = "77+9="
$t.Cell(16, 4).Range.Text = "36-18="
$t.Cell(16, 5).Range.Text = "95-29="
$t.Cell(17, 1).Range.Text = "32+59="
$t.Cell(17, 2).Range.Text = "63-5="
$t.Cell(17, 3).Range.Text = "9+53="
$t.Cell(17, 4).Range.Text = "93-74="
$t.Cell(17, 5).Range.Text = "64-59="
$t.Cell(18, 1).Range.Text = "97-79="
$t.Cell(18, 2).Range.Text = "29+13="
$t.Cell(18, 3).Range.Text = "26+58="
$t.Cell(18, 4).Range.Text = "90-79="
$t.Cell(18, 5).Range.Text = "63-5="
$t.Cell(19, 1).Range.Text = "4+17="
$t.Cell(19, 2).Range.Text = "50-21="
$t.Cell(19, 3).Range.Text = "25+67="
$t.Cell(19, 4).Range.Text = "70-17="
$t.Cell(19, 5).Range.Text = "91-64="
$t.Cell(20, 1).Range.Text = "63-4="
$t.Cell(20, 2).Range.Text = "19+67="
$t.Cell(20, 3).Range.Text = "74-39="
$t.Cell(20, 4).Range.Text = "72-43="
$t.Cell(20, 5).Range.Text = "8+25="
